$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 29412190
$ws.Range("I28").Value = 50000464
$ws.Range("J28").Value = 372.42856
$ws.Range("K28").Value = 50000464
$ws.Range("L28").Value = 372.42856
$ws.Range("M28").Value = -49999979
$ws.Range("N28").Value = -1342.42856

$ws.Range("H127").Value = 7503.5557
$ws.Range("I127").Value = 11429.333
$ws.Range("J127").Value = 3577.7778
$ws.Range("K127").Value = 34287.999
$ws.Range("L127").Value = 10733.3334
$ws.Range("M127").Value = -29327.999
$ws.Range("N127").Value = -20653.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8122631.5
$ws.Range("I32").Value = 2124472.5
$ws.Range("K32").Value = 2124472.5
$ws.Range("M32").Value = -2124185.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H47").Value = 123000
$ws.Range("J47").Value = 123000
$ws.Range("L47").Value = 123000
$ws.Range("N47").Value = -124040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7170.615
$ws.Range("I86").Value = 10541.385
$ws.Range("J86").Value = 3799.8462
$ws.Range("K86").Value = 10541.385
$ws.Range("L86").Value = 3799.8462
$ws.Range("M86").Value = -9418.385
$ws.Range("N86").Value = -6045.8462

$ws.Range("H89").Value = 7170.615
$ws.Range("I89").Value = 10541.385
$ws.Range("J89").Value = 3799.8462
$ws.Range("K89").Value = 52706.925
$ws.Range("L89").Value = 18999.231
$ws.Range("M89").Value = -47090.925
$ws.Range("N89").Value = -30231.231

$ws.Range("H94").Value = 50009652
$ws.Range("J94").Value = 55566188
$ws.Range("L94").Value = 55566188
$ws.Range("N94").Value = -55567090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2992401
$ws.Range("J5").Value = 2565159.8
$ws.Range("L5").Value = 7695479.399999999
$ws.Range("N5").Value = -7695703.399999999

$ws.Range("H55").Value = 1552.8422
$ws.Range("J55").Value = 1892.8572
$ws.Range("L55").Value = 5678.571599999999
$ws.Range("N55").Value = -6032.571599999999

$ws.Range("H68").Value = 17550.834
$ws.Range("I68").Value = 502
$ws.Range("J68").Value = 20960.6
$ws.Range("K68").Value = 1506
$ws.Range("L68").Value = 62881.8
$ws.Range("M68").Value = -695
$ws.Range("N68").Value = -64503.8

$ws.Range("H71").Value = 17550.834
$ws.Range("I71").Value = 502
$ws.Range("J71").Value = 20960.6
$ws.Range("K71").Value = 4518
$ws.Range("L71").Value = 188645.4
$ws.Range("M71").Value = -462
$ws.Range("N71").Value = -196757.4

$ws.Range("H80").Value = 17890.715
$ws.Range("I80").Value = 8600.666999999999
$ws.Range("J80").Value = 20424.363
$ws.Range("K80").Value = 25802.001
$ws.Range("L80").Value = 61273.08900000001
$ws.Range("M80").Value = -24866.001
$ws.Range("N80").Value = -63145.08900000001

$ws.Range("H83").Value = 17890.715
$ws.Range("I83").Value = 8600.666999999999
$ws.Range("J83").Value = 20424.363
$ws.Range("K83").Value = 77406.003
$ws.Range("L83").Value = 183819.267
$ws.Range("M83").Value = -72726.003
$ws.Range("N83").Value = -193179.267

$ws.Range("H114").Value = 1493.5
$ws.Range("I114").Value = 445.57144
$ws.Range("J114").Value = 2716.0833
$ws.Range("K114").Value = 1336.71432
$ws.Range("L114").Value = 8148.249899999999
$ws.Range("M114").Value = 1917.28568
$ws.Range("N114").Value = -14656.2499

$ws.Range("H135").Value = 2992401
$ws.Range("J135").Value = 2565159.8
$ws.Range("L135").Value = 23086438.2
$ws.Range("N135").Value = -23091508.2

$ws.Range("H137").Value = 2056.92
$ws.Range("I137").Value = 1614.5333
$ws.Range("J137").Value = 2720.5
$ws.Range("K137").Value = 4843.5999
$ws.Range("L137").Value = 8161.5
$ws.Range("M137").Value = 256.4000999999998
$ws.Range("N137").Value = -18361.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 25000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 25000
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = 25000
$ws.Range("N15").Value = -25576
$ws.Range("L15").ClearContents()

$ws.Range("H81").Value = 25000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 25000
$ws.Range("K81").Value = 0
$ws.Range("M81").Value = 25000
$ws.Range("N81").Value = -26996
$ws.Range("L81").ClearContents()

$ws.Range("H84").Value = 25000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 25000
$ws.Range("K84").Value = 0
$ws.Range("M84").Value = 75000
$ws.Range("N84").Value = -84984
$ws.Range("L84").ClearContents()

$ws.Range("H102").Value = 5910.1577
$ws.Range("I102").Value = 6399.647
$ws.Range("J102").Value = 1749.5
$ws.Range("K102").Value = 6399.647
$ws.Range("L102").Value = 1749.5
$ws.Range("M102").Value = -4777.647
$ws.Range("N102").Value = -4993.5

$ws.Range("H113").Value = 35598.625
$ws.Range("I113").Value = 6397.8
$ws.Range("J113").Value = 84266.664
$ws.Range("K113").Value = 6397.8
$ws.Range("L113").Value = 84266.664
$ws.Range("M113").Value = -4227.8
$ws.Range("N113").Value = -88606.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 628.6923
$ws.Range("J46").Value = 481.33334
$ws.Range("L46").Value = 481.33334
$ws.Range("N46").Value = -857.33334

$ws.Range("H74").Value = 1000000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 1000000
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = 1000000
$ws.Range("N74").Value = -1001996
$ws.Range("L74").ClearContents()

$ws.Range("H75").Value = 1000000
$ws.Range("J75").Value = 1000000
$ws.Range("L75").Value = 1000000
$ws.Range("N75").Value = -1001872

$ws.Range("H76").Value = 258749.5
$ws.Range("J76").Value = 343332.66
$ws.Range("L76").Value = 343332.66
$ws.Range("N76").Value = -344008.66

$ws.Range("H77").Value = 1000000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1000000
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = 3000000
$ws.Range("N77").Value = -3009984
$ws.Range("L77").ClearContents()

$ws.Range("H78").Value = 1000000
$ws.Range("J78").Value = 1000000
$ws.Range("L78").Value = 3000000
$ws.Range("N78").Value = -3009360

$ws.Range("H79").Value = 258749.5
$ws.Range("J79").Value = 343332.66
$ws.Range("L79").Value = 343332.66
$ws.Range("N79").Value = -345672.66

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16600
$ws.Range("I81").Value = 2100.25
$ws.Range("J81").Value = 20224.938
$ws.Range("K81").Value = 4200.5
$ws.Range("L81").Value = 40449.876
$ws.Range("M81").Value = -3139.5
$ws.Range("N81").Value = -42571.876

$ws.Range("H84").Value = 16600
$ws.Range("I84").Value = 2100.25
$ws.Range("J84").Value = 20224.938
$ws.Range("K84").Value = 21002.5
$ws.Range("L84").Value = 202249.38
$ws.Range("M84").Value = -15698.5
$ws.Range("N84").Value = -212857.38

